$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.321.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.339.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.56%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.73%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.329.34"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.586"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.879.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "599.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.412.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.353.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.903"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.70%  "
$ws.Range("E23").Value = "  -8.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.70%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -16.16%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "559.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.87%  "
$ws.Range("B35").Value = "Cosmos"
$ws.Range("C35").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.840.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "57.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0706"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -13.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.126"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +18.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.342"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "31.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0412"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.129"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
